$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16: GFG-style entry with a numeric question identifier (no GFG label in col A)
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = "Add Two Numbers"

# Row 17: GFG entry (col A labeled "GFG", unstyled like other GFG rows)
$ws.Range("A17").Value = "GFG"
$ws.Range("B17").Value = "Add Two Numbers(Non-reversed)-variation of Leetcode problem-2"

# Apply the same left-aligned, wrap-text formatting used by the other recent rows
# (column B of every data row, plus column A only where it holds a number).
$ws.Range("A16").HorizontalAlignment = -4131
$ws.Range("A16").WrapText = $true

$ws.Range("B16").HorizontalAlignment = -4131
$ws.Range("B16").WrapText = $true

$ws.Range("B17").HorizontalAlignment = -4131
$ws.Range("B17").WrapText = $true

$ws.Range("B18").Select()
